$d = $word.ActiveDocument

# Locate the three consecutive "Co so du lieu" list-item paragraphs that
# together spell out "User (id, name, pass)", "Bai do (...)", and the
# " Info (...)" run, and collapse them into a single, non-list paragraph
# that only contains a single space.
$target = $null
$n = $d.Paragraphs.Count
for ($i = 1; $i -le ($n - 2); $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("User (id, name, pass)")) {
        $next1 = $d.Paragraphs.Item($i + 1).Range.Text
        $next2 = $d.Paragraphs.Item($i + 2).Range.Text
        if ($next1.StartsWith("B") -and $next2.Contains("Info")) {
            $target = $i
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not find the 'User (id, name, pass)' paragraph block"
}

$pStart = $d.Paragraphs.Item($target)
$pEnd = $d.Paragraphs.Item($target + 2)

$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr>' +
    '<w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="28"/>' +
    '<w:szCs w:val="28"/>' +
    '<w:lang w:val="en-GB"/>' +
    '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
    '<w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="28"/>' +
    '<w:szCs w:val="28"/>' +
    '<w:lang w:val="en-GB"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve"> </w:t>' +
    '</w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$r.InsertXML($xml)
